# Order of operations matters: the workbook's shared-string table records
# each distinct text value the first time it is written, and the target
# file expects the new strings in this exact order: firstname, lastname,
# zip, QE-28, test, QE-29. So the raw values below are written first (in
# that order) before any style/formatting passes run.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DATA")

# --- New header cells (row 1): registers "firstname", "lastname", "zip" ---
$ws.Range("G1").Value = "firstname"
$ws.Range("H1").Value = "lastname"
$ws.Range("I1").Value = "zip"

# --- New rows' A column values: registers "QE-28" then (later) "QE-29" ---
$ws.Range("A6").Value = "QE-28"
$ws.Range("A7").Value = "QE-28"

# --- Remaining "test" values (registers "test" after "QE-28") ---
$ws.Range("G3").Value = "test"
$ws.Range("H3").Value = "test"
$ws.Range("G4").Value = "test"
$ws.Range("H4").Value = "test"
$ws.Range("G5").Value = "test"
$ws.Range("H5").Value = "test"
$ws.Range("G6").Value = "test"
$ws.Range("H6").Value = "test"
$ws.Range("H7").Value = "test"

# --- Last new row's A column value: registers "QE-29" ---
$ws.Range("A8").Value = "QE-29"

# --- Numeric cells ---
$ws.Range("I3").Value = 111
$ws.Range("I4").Value = 111
$ws.Range("I5").Value = 111
$ws.Range("I7").Value = 111

# --- Remaining plain text cells for the new rows (columns B-F) ---
$ws.Range("B6").Value = "yes"
$ws.Range("C6").Value = "chrome"
$ws.Range("D6").Value = "standard_user"
$ws.Range("E6").Value = "secret_sauce"
$ws.Range("F6").Value = "Sauce Labs Backpack"

$ws.Range("B7").Value = "yes"
$ws.Range("C7").Value = "chrome"
$ws.Range("D7").Value = "standard_user"
$ws.Range("E7").Value = "secret_sauce"
$ws.Range("F7").Value = "Sauce Labs Backpack, Test.allTheThings() T-Shirt (Red)"

$ws.Range("B8").Value = "yes"
$ws.Range("C8").Value = "chrome"
$ws.Range("D8").Value = "standard_user"
$ws.Range("E8").Value = "secret_sauce"

# --- Empty "quote-prefixed" cells (typed as a lone apostrophe => blank,
#     text-flagged cell, matching the existing F2 style source) ---
$ws.Range("G2").Value = "'"
$ws.Range("H2").Value = "'"
$ws.Range("I2").Value = "'"
$ws.Range("I6").Value = "'"
$ws.Range("G7").Value = "'"
$ws.Range("F8").Value = "'"
$ws.Range("G8").Value = "'"
$ws.Range("H8").Value = "'"
$ws.Range("I8").Value = "'"

# --- Apply the quote-prefix style (copied from the existing F2 cell) to
#     every new cell that needs it ---
$ws.Range("F2").Copy()
$ws.Range("G2:I2").PasteSpecial(-4122)
$ws.Range("G3:I3").PasteSpecial(-4122)
$ws.Range("G4:I4").PasteSpecial(-4122)
$ws.Range("G5:I5").PasteSpecial(-4122)
$ws.Range("I6").PasteSpecial(-4122)
$ws.Range("G7").PasteSpecial(-4122)
$ws.Range("F8:I8").PasteSpecial(-4122)

# --- Update the selected cell to match the saved view state ---
$ws.Range("E25").Select()
